$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (rows 2-15) holds a date serial number (45171 -> 2023-09-02)
# that needs to be bumped by one day to 45172 (2023-09-03).
for ($row = 2; $row -le 15; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45171) {
        $cell.Value = 45172
    }
}
